# feat: add 2022-Q1 data
#
# - Insert a new sheet "2022-Q1" right before the "总计" (totals) sheet,
#   populated with per-fund holding data (same layout as the other
#   quarterly sheets: "2021-Q1" / "2021-Q2" / "2021-Q3").
# - Update the "总计" sheet: add a new first data row for 2022-Q1
#   (2 funds, 0.6 billion yuan) and push the existing quarters down.

$wb = $excel.ActiveWorkbook

# Template sheet to copy the header/index-column formatting from for the
# new quarterly sheet (all quarterly sheets share the same look).
$template = $wb.Worksheets.Item("2021-Q2")

# ---------------------------------------------------------------------
# 1) Create the "2022-Q1" worksheet, positioned right before "总计".
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$newSheet = $wb.Worksheets.Add($totalSheet)
$newSheet.Name = "2022-Q1"

# NOTE: after Add(), the *old* $totalSheet handle now resolves to the
# newly-inserted sheet (positional aliasing), so re-fetch "总计" by name
# to get a handle that actually points at the totals sheet again.
$totalSheet = $wb.Worksheets.Item("总计")

# Helper to write a plain-text value into a cell without Excel coercing
# numeric-looking strings (e.g. "18.44", "001186") into numbers - the
# leading apostrophe forces text entry, then formatting is cleared so no
# stray "quote prefix" style lingers on the cell (these body cells use
# the sheet's default formatting, same as the source data).
function Set-TextCell($ws, $addr, $text) {
    $ws.Range($addr).Value = "'" + $text
    $ws.Range($addr).ClearFormats()
}

# Bring over the bold/centered/bordered look of the header row (B1:H1)
# and the index column (A2, A3) from an existing quarterly sheet.
$template.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)
$template.Range("A2").Copy()
$newSheet.Range("A2:A3").PasteSpecial(-4122)

# Header row
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Row 2 - 富国文体健康股票A
$newSheet.Range("A2").Value = 0
Set-TextCell $newSheet "B2" "001186"
Set-TextCell $newSheet "C2" "富国文体健康股票A"
Set-TextCell $newSheet "D2" "18.44"
Set-TextCell $newSheet "E2" "86.61"
Set-TextCell $newSheet "F2" "3.15"
Set-TextCell $newSheet "G2" "0.5809"
$newSheet.Range("H2").Value = 10

# Row 3 - 富国文体健康股票C
$newSheet.Range("A3").Value = 1
Set-TextCell $newSheet "B3" "011125"
Set-TextCell $newSheet "C3" "富国文体健康股票C"
Set-TextCell $newSheet "D3" "0.71"
Set-TextCell $newSheet "E3" "86.61"
Set-TextCell $newSheet "F3" "3.15"
Set-TextCell $newSheet "G3" "0.0224"
$newSheet.Range("H3").Value = 10

# ---------------------------------------------------------------------
# 2) Update the "总计" sheet with the new 2022-Q1 summary row, shifting
#    the pre-existing quarters down by one row.
# ---------------------------------------------------------------------

# Extend the formatted index column (A2:A4 already carry the bold
# centered/bordered look) down to the new row 5.
$totalSheet.Range("A4").Copy()
$totalSheet.Range("A5").PasteSpecial(-4122)

$totalSheet.Range("A2").Value = 0
Set-TextCell $totalSheet "B2" "2022-Q1"
$totalSheet.Range("C2").Value = 2
$totalSheet.Range("D2").Value = 0.6

$totalSheet.Range("A3").Value = 1
Set-TextCell $totalSheet "B3" "2021-Q3"
$totalSheet.Range("C3").Value = 1
$totalSheet.Range("D3").Value = 0.02

$totalSheet.Range("A4").Value = 2
Set-TextCell $totalSheet "B4" "2021-Q2"
$totalSheet.Range("C4").Value = 3
$totalSheet.Range("D4").Value = 0.03

$totalSheet.Range("A5").Value = 3
Set-TextCell $totalSheet "B5" "2021-Q1"
$totalSheet.Range("C5").Value = 3
$totalSheet.Range("D5").Value = 0.06

# Restore the originally-active sheet/tab (this edit shouldn't change
# which sheet the workbook opens to).
$wb.Worksheets.Item("2021-Q1").Activate()
